$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Create rows 82 and 83 by copying the formatting of row 81 (keeps styles for A:K incl blanks I/J/K) ---
$ws.Range("A81:K81").Copy($ws.Range("A82:K82"))
$ws.Range("A81:K81").Copy($ws.Range("A83:K83"))

# --- Update row 81: existing "Check Out Transaction" test gains many new params, height grows ---
$ws.Range("H81").Value = "coyni.checkout.tests.CheckOutTest,`ntesCheckOutTransaction,`n-pdomain,`n-ppublicKey,`n-psecretKey,`n-pemail,`n-ppassword,`n-pheading,`n-pcode,`n-psuccessContent,`n-pinsufficient,`n-pcvv,`n-pamount,`n-pnameOnCard,`n-pcardNumber,`n-pcardType,`n-pcardExpiry,`n-pcvvNumber,`n-paddressLine1,`n-paddressLine2,`n-pcity,`n-pzipCode,`n-pstate,`n-pcountry,`n-ppreamount`n"
$ws.Rows.Item(81).RowHeight = 390

# --- Row 82: new "invalid PopUp" check-out test ---
$ws.Range("A82").Value = "test checkOut Transaction invalid PopUp"
$ws.Range("H82").Value = "coyni.checkout.tests.CheckOutTest,`ntesCheckOutInvalidData,`n-pdomain,`n-porderId,`n-ppublicKey,`n-psecretKey"
$ws.Rows.Item(82).RowHeight = 90

# --- Row 83: new "Cancel Transaction" check-out test ---
$ws.Range("A83").Value = "test checkOut Cancel Transaction"
$ws.Range("H83").Value = "coyni.checkout.tests.CheckOutTest,`ntesCheckOutCancelTransaction,`n-pdomain,`n-ppublicKey,`n-psecretKey,`n-pemail,`n-ppassword,`n-pheading,`n-pcode,`n-pcancelHeading,`n-pcontent,`n-pcheckOutContent`n"
$ws.Rows.Item(83).RowHeight = 195

# --- Update the sheet view so the newly added rows are visible/selected ---
$ws.Range("B83").Select()
